# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# This workbook is a "Estado de Cuenta" (account statement) listing, per
# worker, the late-payment periods owed. A new period (2509) needs to be
# added for the two workers already present in the table (row 19/20 which
# carry period 2508), and the summary totals (Valor Mora / Cant. Periodos)
# need to be refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The detail table currently ends at row 20 (period 2508, two rows: CC
# 45622751 and CC 1052210104). Row 21 onward is blank until the signature
# block, which currently lives on rows 25-26.
#
# Insert two new blank rows right above the signature block (old row 25)
# so it is pushed down to rows 27-28, matching the rest of the table's
# layout, and leaving rows 21-22 free for the new "2509" period entries.
$ws.Rows("25:26").Insert()

# Duplicate the formatting of the last two existing data rows (19:20, the
# 2508 entries) onto the two freshly inserted rows (21:22) so the new
# period keeps the same borders/shading/number formats as the rest of the
# table.
$ws.Range("B19:J20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4142)
$excel.CutCopyMode = $false

# Row 21: CC 45622751 - LINA MARCELA PEREZ CAICEDO - periodo 2509
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45622751"
$ws.Range("D21").Value = "LINA MARCELA PEREZ CAICEDO"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# Row 22: CC 1052210104 - SERGIO ANDRES MONSALVE RICO - periodo 2509
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1052210104"
$ws.Range("D22").Value = "SERGIO ANDRES MONSALVE RICO"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 120000
$ws.Range("G22").Value = 3000000

# Refresh the "VALOR MORA" total (E11) to include the two new rows and the
# "Cant. Periodos" count (F13), which now spans four distinct periods
# (2506, 2507, 2508, 2509) instead of three.
$ws.Range("E11").Value = 622820
$ws.Range("F13").Value = 4
